$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.219740271568298
$ws.Range("B1").Value = 2.687123775482178
$ws.Range("C1").Value = 4.3643798828125
$ws.Range("D1").Value = 2.145384311676025
$ws.Range("E1").Value = 1.16199791431427
